$wb = $excel.ActiveWorkbook

# --- Step 1: Create the new "2022-Q1" sheet by copying the "2021-Q4" template ---
# (this preserves header text/styles, index-column style, dimension, etc.)
$totalSheetBeforeCopy = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheetBeforeCopy)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Re-fetch "总计" by name: inserting the new sheet shifted sheet positions, so any
# handle captured before the Copy() now points at the wrong sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# The template only had 16 data rows; extend the bold-bordered index-column style
# (column A) down to the extra rows we need before filling them in.
$newSheet.Range("A2").Copy()
$newSheet.Range("A18:A23").PasteSpecial(-4122)

# Force text-typed columns (B,C,D,E,F,G) to stay text even for numeric-looking values
# (fund codes like "006682" must not collapse to 6682, decimals must stay literal strings)
$textRange = $newSheet.Range("B2:G23")
$textRange.NumberFormat = "@"

# --- Step 2: Write the 2022-Q1 fund holding rows ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "006682"
$newSheet.Range("C2").Value = "景顺长城中证500指数增强"
$newSheet.Range("D2").Value = "16.63"
$newSheet.Range("E2").Value = "87.75"
$newSheet.Range("F2").Value = "1.97"
$newSheet.Range("G2").Value = "0.3276"
$newSheet.Range("H2").Value = 8
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "003318"
$newSheet.Range("C3").Value = "景顺长城中证500行业中性低波动指数"
$newSheet.Range("D3").Value = "13.99"
$newSheet.Range("E3").Value = "93.88"
$newSheet.Range("F3").Value = "1.71"
$newSheet.Range("G3").Value = "0.2392"
$newSheet.Range("H3").Value = 3
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "005994"
$newSheet.Range("C4").Value = "国投瑞银中证500指数量化增强A"
$newSheet.Range("D4").Value = "11.53"
$newSheet.Range("E4").Value = "87.00"
$newSheet.Range("F4").Value = "1.47"
$newSheet.Range("G4").Value = "0.1695"
$newSheet.Range("H4").Value = 4
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "000978"
$newSheet.Range("C5").Value = "景顺长城量化精选股票"
$newSheet.Range("D5").Value = "8.51"
$newSheet.Range("E5").Value = "93.86"
$newSheet.Range("F5").Value = "1.99"
$newSheet.Range("G5").Value = "0.1693"
$newSheet.Range("H5").Value = 4
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "001050"
$newSheet.Range("C6").Value = "汇添富成长多因子量化策略股票"
$newSheet.Range("D6").Value = "11.48"
$newSheet.Range("E6").Value = "92.68"
$newSheet.Range("F6").Value = "0.71"
$newSheet.Range("G6").Value = "0.0815"
$newSheet.Range("H6").Value = 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "008851"
$newSheet.Range("C7").Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$newSheet.Range("D7").Value = "5.05"
$newSheet.Range("E7").Value = "74.55"
$newSheet.Range("F7").Value = "1.57"
$newSheet.Range("G7").Value = "0.0793"
$newSheet.Range("H7").Value = 7
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "009992"
$newSheet.Range("C8").Value = "景顺长城量化成长演化混合"
$newSheet.Range("D8").Value = "2.65"
$newSheet.Range("E8").Value = "92.88"
$newSheet.Range("F8").Value = "2.55"
$newSheet.Range("G8").Value = "0.0676"
$newSheet.Range("H8").Value = 5
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "006511"
$newSheet.Range("C9").Value = "博道卓远混合A"
$newSheet.Range("D9").Value = "2.20"
$newSheet.Range("E9").Value = "82.79"
$newSheet.Range("F9").Value = "2.57"
$newSheet.Range("G9").Value = "0.0565"
$newSheet.Range("H9").Value = 10
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "007089"
$newSheet.Range("C10").Value = "国投瑞银中证500指数量化增强C"
$newSheet.Range("D10").Value = "3.82"
$newSheet.Range("E10").Value = "87.00"
$newSheet.Range("F10").Value = "1.47"
$newSheet.Range("G10").Value = "0.0562"
$newSheet.Range("H10").Value = 4
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "003016"
$newSheet.Range("C11").Value = "中金中证500指数增强A"
$newSheet.Range("D11").Value = "4.71"
$newSheet.Range("E11").Value = "93.78"
$newSheet.Range("F11").Value = "1.12"
$newSheet.Range("G11").Value = "0.0528"
$newSheet.Range("H11").Value = 9
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "007825"
$newSheet.Range("C12").Value = "博道志远混合A"
$newSheet.Range("D12").Value = "1.67"
$newSheet.Range("E12").Value = "82.63"
$newSheet.Range("F12").Value = "2.55"
$newSheet.Range("G12").Value = "0.0426"
$newSheet.Range("H12").Value = 10
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "005258"
$newSheet.Range("C13").Value = "景顺长城量化平衡灵活配置混合"
$newSheet.Range("D13").Value = "2.39"
$newSheet.Range("E13").Value = "90.00"
$newSheet.Range("F13").Value = "1.60"
$newSheet.Range("G13").Value = "0.0382"
$newSheet.Range("H13").Value = 8
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "512260"
$newSheet.Range("C14").Value = "华安中证500行业中性低波动ETF"
$newSheet.Range("D14").Value = "1.17"
$newSheet.Range("E14").Value = "96.94"
$newSheet.Range("F14").Value = "1.77"
$newSheet.Range("G14").Value = "0.0207"
$newSheet.Range("H14").Value = 3
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "011731"
$newSheet.Range("C15").Value = "国投瑞银安睿混合A"
$newSheet.Range("D15").Value = "2.58"
$newSheet.Range("E15").Value = "43.48"
$newSheet.Range("F15").Value = "0.73"
$newSheet.Range("G15").Value = "0.0188"
$newSheet.Range("H15").Value = 5
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "007826"
$newSheet.Range("C16").Value = "博道志远混合C"
$newSheet.Range("D16").Value = "0.73"
$newSheet.Range("E16").Value = "82.63"
$newSheet.Range("F16").Value = "2.55"
$newSheet.Range("G16").Value = "0.0186"
$newSheet.Range("H16").Value = 10
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "003578"
$newSheet.Range("C17").Value = "中金中证500指数增强C"
$newSheet.Range("D17").Value = "1.44"
$newSheet.Range("E17").Value = "93.78"
$newSheet.Range("F17").Value = "1.12"
$newSheet.Range("G17").Value = "0.0161"
$newSheet.Range("H17").Value = 9
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "006729"
$newSheet.Range("C18").Value = "万家中证500指数增强A"
$newSheet.Range("D18").Value = "1.04"
$newSheet.Range("E18").Value = "93.64"
$newSheet.Range("F18").Value = "1.27"
$newSheet.Range("G18").Value = "0.0132"
$newSheet.Range("H18").Value = 6
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "004192"
$newSheet.Range("C19").Value = "招商中证500指数增强A"
$newSheet.Range("D19").Value = "0.96"
$newSheet.Range("E19").Value = "94.32"
$newSheet.Range("F19").Value = "1.14"
$newSheet.Range("G19").Value = "0.0109"
$newSheet.Range("H19").Value = 5
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "006730"
$newSheet.Range("C20").Value = "万家中证500指数增强C"
$newSheet.Range("D20").Value = "0.61"
$newSheet.Range("E20").Value = "93.64"
$newSheet.Range("F20").Value = "1.27"
$newSheet.Range("G20").Value = "0.0077"
$newSheet.Range("H20").Value = 6
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "011732"
$newSheet.Range("C21").Value = "国投瑞银安睿混合C"
$newSheet.Range("D21").Value = "0.95"
$newSheet.Range("E21").Value = "43.48"
$newSheet.Range("F21").Value = "0.73"
$newSheet.Range("G21").Value = "0.0069"
$newSheet.Range("H21").Value = 5
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22").Value = "004193"
$newSheet.Range("C22").Value = "招商中证500指数增强C"
$newSheet.Range("D22").Value = "0.42"
$newSheet.Range("E22").Value = "94.32"
$newSheet.Range("F22").Value = "1.14"
$newSheet.Range("G22").Value = "0.0048"
$newSheet.Range("H22").Value = 5
$newSheet.Range("A23").Value = 21
$newSheet.Range("B23").Value = "006512"
$newSheet.Range("C23").Value = "博道卓远混合C"
$newSheet.Range("D23").Value = "0.10"
$newSheet.Range("E23").Value = "82.79"
$newSheet.Range("F23").Value = "2.57"
$newSheet.Range("G23").Value = "0.0026"
$newSheet.Range("H23").Value = 10

# Clear the temporary text numberformat back to default (keeps cell *type* as Text,
# just removes the explicit @ style so it matches the rest of the workbook)
$textRange.Style = "Normal"

# --- Step 3: Update the "总计" sheet: insert a new top row for 2022-Q1, shift the rest down ---
$totalSheet.Rows.Item(2).Insert()

# Excel inherits a blended style from the header row on insert; re-stamp row 2 with the
# same plain data-row formatting used by every other row (copy format from row 3, which
# is the original "2021-Q4" row that just got pushed down).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 22
$totalSheet.Range("D2").Value = 1.5
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 16
$totalSheet.Range("D3").Value = 1.49
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 8
$totalSheet.Range("D4").Value = 0.63
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 14
$totalSheet.Range("D5").Value = 0.66
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 3
$totalSheet.Range("D6").Value = 0.72
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 15
$totalSheet.Range("D7").Value = 3.38

# Restore the originally active sheet/tab (the copy/rename dance above left the new
# sheet focused) so the workbook-level view state is otherwise untouched.
$wb.Worksheets.Item(1).Activate()

Write-Output "done"
